$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the year value in B1 with the header label "Valor"
$ws.Range("B1").Value = "Valor"

# Set explicit best-fit widths for columns B and C
$ws.Columns.Item(2).ColumnWidth = 7.0
$ws.Columns.Item(3).ColumnWidth = 6.833333333333334

# Move the active selection to B3
$ws.Range("B3").Select() | Out-Null
